# Fruta / hortaliza, semanal
# Insert 4 new weekly data rows (new date 44529) right after the existing
# row 50 ("Femacal de La Calera" / Damasco data), pushing the previous
# rows 51-68 down to rows 55-72, then populate the newly inserted rows
# 51-54 with the new week's observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 51, shifting rows 51-68 down
# to become rows 55-72.
$ws.Rows.Item(51).Resize(4).Insert()

$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44529, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Castle Brite", "Especial", 56, 20000, 20000, 20000, "`$/caja 15 kilos granel", "Provincia de San Felipe de Aconcagua", 1333, 15),
    @(3, "Femacal de La Calera", "Coquimbo", 44529, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Castle Brite", "Primera", 50, 18000, 18000, 18000, "`$/caja 15 kilos granel", "Provincia de San Felipe de Aconcagua", 1200, 15),
    @(3, "Femacal de La Calera", "Coquimbo", 44529, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Modesto", "Especial", 56, 14000, 14000, 14000, "`$/bandeja 10 kilos", "Provincia de San Felipe de Aconcagua", 1400, 10),
    @(3, "Femacal de La Calera", "Coquimbo", 44529, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Modesto", "Primera", 50, 12000, 12000, 12000, "`$/bandeja 10 kilos", "Provincia de San Felipe de Aconcagua", 1200, 10)
)

$startRow = 51
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowData = $newRows[$i]
    $r = $startRow + $i
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowData[$col - 1]
    }
}

# Apply the same date number format used by column D elsewhere to the
# newly written date cells.
$ws.Range("D51:D54").NumberFormat = $ws.Range("D50").NumberFormat
